# Add the new weekly data column (week of 11_05_2021) as column AB,
# following the same pattern as the existing weekly columns (B..AA).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for the new week
$ws.Range("AB1").Value = "11_05_2021"

# Age-group death counts for the new week
$ws.Range("AB2").Value = 1
$ws.Range("AB3").Value = 0
$ws.Range("AB4").Value = 0
$ws.Range("AB5").Value = 7
$ws.Range("AB6").Value = 8
$ws.Range("AB7").Value = 62
$ws.Range("AB8").Value = 208
$ws.Range("AB9").Value = 661
$ws.Range("AB10").Value = 991
$ws.Range("AB11").Value = 561

# Total row: same SUM pattern used by the other weekly columns
$ws.Range("AB12").Formula = "=SUM(AB2:AB11)"

# Mirror the scrolled/selected view from the edited workbook as closely
# as possible.
$ws.Range("T1").Select()
$ws.Range("AB14").Select()
